$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Andre Russell" values in the source data end with a non-breaking space
# (U+00A0), matching the existing F2/F3 cells in the workbook.
$nbsp = [char]0x00A0
$player = "Andre Russell$nbsp"

# Insert a brand-new row at position 2 (Oct 12 2020 vs RCB), pushing the
# existing Oct 16 2020 / Oct 7 2020 rows down to rows 3 and 4.
$ws.Rows("2:2").Insert()

# All cells in this sheet are stored as literal text (t="str"), even the
# numeric-looking ones (runs/balls/4s/6s/sr) - that's why the sheet also
# carries a numberStoredAsText ignoredError. Force text formatting on the
# whole target block before writing so numeric-looking strings like "16"
# or "160.00" don't get auto-coerced into real numbers, then clear the
# formatting afterwards so no stray style index is left on the cells
# (matching the original, which has no explicit cell styles).
$target = $ws.Range("A2:K10")
$target.NumberFormat = "@"

function Set-RowValues($rowIndex, $values) {
    $arr = New-Object 'object[,]' 1,11
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0, $i] = $values[$i]
    }
    $ws.Range("A" + $rowIndex + ":K" + $rowIndex).Value = $arr
}

Set-RowValues 2  @(" Oct 12 2020", " Sharjah", "RCB won by 82 runs", "Kolkata Knight Riders", "Royal Challengers Bangalore", $player, "16", "10", "2", "1", "160.00")
Set-RowValues 3  @(" Oct 16 2020", " Abu Dhabi", "Mumbai won by 8 wickets (with 19 balls remaining)", "Kolkata Knight Riders", "Mumbai Indians", $player, "12", "9", "1", "1", "133.33")
Set-RowValues 4  @(" Oct 3 2020", " Sharjah", "Capitals won by 18 runs", "Kolkata Knight Riders", "Delhi Capitals", $player, "13", "8", "1", "1", "162.50")
Set-RowValues 5  @(" Oct 18 2020", " Abu Dhabi", "Match tied (KKR won the one-over eliminator)", "Kolkata Knight Riders", "Sunrisers Hyderabad", $player, "9", "11", "1", "0", "81.81")
Set-RowValues 6  @(" Oct 7 2020", " Abu Dhabi", "KKR won by 10 runs", "Kolkata Knight Riders", "Chennai Super Kings", $player, "2", "4", "0", "0", "50.00")
Set-RowValues 7  @(" Oct 10 2020", " Abu Dhabi", "KKR won by 2 runs", "Kolkata Knight Riders", "Kings XI Punjab", $player, "5", "3", "1", "0", "166.66")
Set-RowValues 8  @(" Nov 1 2020", " Dubai (DSC)", "KKR won by 60 runs", "Kolkata Knight Riders", "Rajasthan Royals", $player, "25", "11", "1", "3", "227.27")
Set-RowValues 9  @(" Sep 30 2020", " Dubai (DSC)", "KKR won by 37 runs", "Kolkata Knight Riders", "Rajasthan Royals", $player, "24", "14", "0", "3", "171.42")
Set-RowValues 10 @(" Sep 23 2020", " Abu Dhabi", "Mumbai won by 49 runs", "Kolkata Knight Riders", "Mumbai Indians", $player, "11", "11", "2", "0", "100.00")

$target.ClearFormats()
